$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.392.10'
$ws.Range('E2').Value = '  -0.83%  '
$ws.Range('D3').Value = '2.328.21'
$ws.Range('E3').Value = '  -0.75%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '511.63'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.28%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.82'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.78%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -0.72%  '
$ws.Range('E9').Value = '  -3.23%  '
$ws.Range('E10').Value = '  -0.55%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.28'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.338'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.45%  '
$ws.Range('D13').Value = '2.745.70'
$ws.Range('E13').Value = '  -0.55%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.50'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.56%  '
$ws.Range('D15').Value = '56.397.45'
$ws.Range('E15').Value = '  -0.65%  '
$ws.Range('E16').Value = '  -1.40%  '
$ws.Range('D17').Value = '2.330.44'
$ws.Range('E17').Value = '  -1.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.43'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '324.01'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.13'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -2.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.65'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '61.79'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.76%  '
$ws.Range('E24').Value = '  +11.62%  '
$ws.Range('E25').Value = '  +1.11%  '
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('E27').Value = '  +4.98%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '167.72'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.55%  '
$ws.Range('E29').Value = '  -0.60%  '
$ws.Range('D30').Value = '0.0₃0716'
$ws.Range('E30').Value = '  -3.32%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.09'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.29'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.38%  '
$ws.Range('E34').Value = '  -0.18%  '
$ws.Range('E35').Value = '  +1.08%  '
$ws.Range('E36').Value = '  -1.32%  '
$ws.Range('E37').Value = '  -4.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '38.42'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.77%  '
$ws.Range('E39').Value = '  +0.91%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '150.77'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +10.15%  '
$ws.Range('E41').Value = '  -0.92%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.57'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '277.63'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.39%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.01'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.40%  '
$ws.Range('E45').Value = '  -1.01%  '
$ws.Range('E46').Value = '  -1.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.556'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '17.98'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +5.08%  '
$ws.Range('E49').Value = '  -1.52%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.378'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.13'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +2.05%  '
